$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("RUNMANAGER")
$ws2 = $wb.Worksheets.Item("DATA")

# --- Sheet2 (DATA) content updates -----------------------------------
# Row1 header: browser / username / password columns
$ws2.Range("C1").Value = "browser"
$ws2.Range("D1").Value = "username"
$ws2.Range("E1").Value = "password"

# Row2
$ws2.Range("B2").Value = "yes"
$ws2.Range("C2").Value = "chrome"
$ws2.Range("D2").Value = "spcbtest"
$ws2.Range("E2").Value = "Asdf@123"

# Row3
$ws2.Range("A3").Value = "newTest"
$ws2.Range("B3").Value = "yes"
$ws2.Range("C3").Value = "chrome"
$ws2.Range("D3").Value = "spcb"
$ws2.Range("E3").Value = "Asdf@123"

# Rows 4-6: clear all existing data (keep the E column's hyperlink style only)
$ws2.Range("A4:D6").ClearContents()
$ws2.Range("E4:E6").ClearContents()

# Remove hyperlinks attached to E4, E5, E6 (keep E2 / E3).
# Re-scan the live collection for each address (in reverse order) instead of
# caching hyperlink references up-front, since deleting one shifts/ invalidates
# the others' cached identities.
function Remove-HyperlinkAt($ws, $addr) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.Delete()
            return
        }
    }
}
Remove-HyperlinkAt $ws2 "`$E`$6"
Remove-HyperlinkAt $ws2 "`$E`$5"
Remove-HyperlinkAt $ws2 "`$E`$4"

# Add row 7, matching the hyperlink-style (but empty) cell formatting used by E4:E6
$ws2.Range("E7").Style = $ws2.Range("E6").Style

# --- Sheet1 (RUNMANAGER) content updates ------------------------------
# Row3 changes
$ws1.Range("B3").Value = "To check this test runs"
$ws1.Range("C3").Value = "yes"
$ws1.Range("D3").Value = "'1"

# --- Selections / active sheet ---------------------------------------
# Set sheet2's selection first (this also makes it momentarily active),
# then select on sheet1 last so RUNMANAGER ends up the active/selected tab.
$ws2.Range("C20").Select()
$ws1.Range("C10").Select()
$ws1.Activate()

# --- Rename sheet2 ------------------------------------------------------
$ws2.Name = "RETAIL_DATA"
